$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dmanisi")

# Copy formatting from column J (2022) into the new column K (2023)
# so the new cells reuse the same styles as the rest of the table.
$ws.Range("J3:J6").Copy() | Out-Null
$ws.Range("K3:K6").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Fill in the new 2023 data column
$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 1439.3
$ws.Range("K5").Value = 748
$ws.Range("K6").Value = 1929.3

$excel.CutCopyMode = 0
